$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/25/2025  Through  8/31/2025"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -26.086956521739
$ws.Range("L16").Value = -27.659574468085
$ws.Range("M16").Value = -44.262295081967
$ws.Range("N16").Value = -84.888888888888
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -85.714285714285
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -57.142857142857
$ws.Range("I17").Value = 61
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = -26.506024096385
$ws.Range("L17").Value = -17.567567567567
$ws.Range("M17").Value = -14.084507042253
$ws.Range("N17").Value = -71.090047393364
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 63
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = 3.278688524590
$ws.Range("L18").Value = -8.695652173913
$ws.Range("M18").Value = -13.698630136986
$ws.Range("N18").Value = -77.738515901060
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = -6.666666666666
$ws.Range("I19").Value = 130
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = -7.142857142857
$ws.Range("L19").Value = 19.266055045871
$ws.Range("M19").Value = -12.751677852349
$ws.Range("N19").Value = -12.162162162162
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = -21.621621621621
$ws.Range("L20").Value = -35.555555555555
$ws.Range("M20").Value = -23.684210526315
$ws.Range("N20").Value = -89.138576779026
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 35
$ws.Range("G21").Value = 59
$ws.Range("H21").Value = -40.677966101694
$ws.Range("I21").Value = 319
$ws.Range("J21").Value = 370
$ws.Range("K21").Value = -13.783783783783
$ws.Range("L21").Value = -8.595988538681
$ws.Range("M21").Value = -19.240506329113
$ws.Range("N21").Value = -72.212543554007
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 150
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = -54.545454545454
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = -9.090909090909
$ws.Range("I23").Value = 65
$ws.Range("J23").Value = 75
$ws.Range("K23").Value = -13.333333333333
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = 16.071428571428
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 55.555555555555
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 6.521739130434
$ws.Range("I24").Value = 324
$ws.Range("J24").Value = 351
$ws.Range("K24").Value = -7.692307692307
$ws.Range("L24").Value = -20.393120393120
$ws.Range("M24").Value = 8.361204013377
$ws.Range("C25").Value = 3
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -21.052631578947
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 182
$ws.Range("K25").Value = -29.670329670329
$ws.Range("L25").Value = -39.336492890995
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 12
$ws.Range("H26").Value = 58.333333333333
$ws.Range("I26").Value = 129
$ws.Range("J26").Value = 103
$ws.Range("K26").Value = 25.242718446601
$ws.Range("L26").Value = 19.444444444444
$ws.Range("M26").Value = -31.382978723404
$ws.Range("G28").Value = 1
$ws.Range("N29").Value = -92.857142857142
$ws.Range("N30").Value = -94.117647058823
$ws.Range("K33").Value = 0

# --- Cells that switch from the "0"/"***.*" text placeholder to a real number ---
# Set the numeric value, then copy the number formatting from a stable same-style
# numeric cell (D16, style 15) so the cell style matches the target exactly.
$ws.Range("C16").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C20").Value = 2
$ws.Range("D16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C22").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("D16").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("C33").Value = 1
$ws.Range("D16").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("F33").Value = 1
$ws.Range("D16").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("I33").Value = 1
$ws.Range("D16").Copy()
$ws.Range("I33").PasteSpecial(-4122)

# --- Cells that switch from a real number to the "0"/"***.*" text placeholder ---
# Force text typing with a leading apostrophe, then copy formatting from a stable
# same-style placeholder cell (D22 = "0", E22 = "***.*") so the style matches exactly.
$ws.Range("D28").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("E22").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$excel.CutCopyMode = 0
